$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.846.21"
$ws.Range("E2").Value = "  -1.12%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.562.61"
$ws.Range("E3").Value = "  +0.23%  "
$ws.Range("E4").Value = "  -0.14%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "205.77"
$ws.Range("E5").Value = "  -0.29%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.486"
$ws.Range("E6").Value = "  -1.56%  "
$ws.Range("E7").Value = "  -0.13%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "21.78"
$ws.Range("E8").Value = "  -1.86%  "
$ws.Range("E9").Value = "  -0.50%  "
$ws.Range("E10").Value = "  -1.37%  "
$ws.Range("E11").Value = "  +0.25%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.783.72"
$ws.Range("E12").Value = "  +0.17%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.565.80"
$ws.Range("E13").Value = "  +0.48%  "
$ws.Range("E14").Value = "  -1.32%  "
$ws.Range("E15").Value = "  -0.35%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "26.857.17"
$ws.Range("E16").Value = "  -1.07%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "61.24"
$ws.Range("E17").Value = "  -2.41%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "215.01"
$ws.Range("E18").Value = "  +0.63%  "
$ws.Range("E19").Value = "  +1.82%  "
$ws.Range("E20").Value = "  -0.93%  "
$ws.Range("E21").Value = "  -0.12%  "
$ws.Range("E22").Value = "  +0.36%  "
$ws.Range("E23").Value = "  -1.60%  "
$ws.Range("E24").Value = "  +1.39%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "153.73"
$ws.Range("E25").Value = "  +1.15%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.71"
$ws.Range("E26").Value = "  +1.66%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "14.91"
$ws.Range("E27").Value = "  +0.28%  "
$ws.Range("E28").Value = "  -0.14%  "
$ws.Range("E29").Value = "  -1.31%  "
$ws.Range("E30").Value = "  +0.87%  "
$ws.Range("E31").Value = "  -3.57%  "
$ws.Range("E32").Value = "  -0.11%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.401.47"
$ws.Range("E33").Value = "  +1.54%  "
$ws.Range("E34").Value = "  -0.55%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.53"
$ws.Range("E35").Value = "  -0.98%  "
$ws.Range("E36").Value = "  -0.50%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.913"
$ws.Range("E37").Value = "  -3.57%  "
$ws.Range("E38").Value = "  -0.91%  "
$ws.Range("E39").Value = "  +2.11%  "
$ws.Range("E40").Value = "  -0.24%  "
$ws.Range("E41").Value = "  -0.14%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.995"
$ws.Range("E42").Value = "  +1.25%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.55"
$ws.Range("E43").Value = "  +6.48%  "
$ws.Range("E44").Value = "  -0.22%  "
$ws.Range("E45").Value = "  +1.11%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "63.32"
$ws.Range("E46").Value = "  +0.23%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.697.66"
$ws.Range("E47").Value = "  +0.29%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "86.48"
$ws.Range("E48").Value = "  +1.23%  "
$ws.Range("E49").Value = "  +2.13%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0₇0976"
$ws.Range("E50").Value = "  -0.83%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0948"
$ws.Range("E51").Value = "  +0.89%  "
